# Workbook "ComponentesCamaraUVO" was re-uploaded with updated material
# callouts: the steel grade used for the housing/lid/connector/side-door
# parts changed from SAE1020 to SAE1010 (and a typo "Cielab" -> "Cienlab"
# was fixed on the side door description).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B3").Value = "Casco inferior (SAE1010) (Cienlab)"
$ws.Range("B4").Value = "Casco superior (SAE1010) (Cienlab)"
$ws.Range("B5").Value = "Conexão entre cascos (SAE1010) (Cienlab)"
$ws.Range("B6").Value = "Porta lateral (SAE1010)(Cienlab)"

# Mirror the cursor/selection position recorded in the saved file.
$ws.Range("K14").Select()
